# "Start to delete GlobaleVars" — progress update on the To-Do list:
#   - B5 "Supprimer les variables globales" moves from 0% to 70% complete
#   - B8 "Migrer vers eclipse" is marked fully done (100%)
#   - Active selection ends up on F6 (next row the user was about to touch)
#   - DATA sheet gets an explicit width on column A

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Liste de tâches")

# "Supprimer les variables globales" -> 70% complete
$ws1.Range("F5").Value = 0.7

# "Migrer vers eclipse" -> 100% complete
$ws1.Range("F8").Value = 1

# Leave the cursor on F6, matching where the user navigated to next
$ws1.Range("F6").Select() | Out-Null

# DATA sheet: column A (priority list) gets an explicit custom width
$ws2 = $wb.Worksheets.Item("DATA")
$ws2.Columns.Item(1).ColumnWidth = 11.855
